$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'27.997.31"
$ws.Range('E2').Value = '  -2.24%  '

# Row 3
$ws.Range('D3').Value = "'1.828.85"
$ws.Range('E3').Value = '  -1.22%  '

# Row 4
$ws.Range('E4').Value = '  -0.13%  '

# Row 5
$ws.Range('D5').Value = "'326.17"
$ws.Range('E5').Value = '  -2.97%  '

# Row 6
$ws.Range('E6').Value = '  -0.18%  '

# Row 7
$ws.Range('D7').Value = "'0.4610"
$ws.Range('E7').Value = '  -0.90%  '

# Row 8
$ws.Range('D8').Value = "'0.3859"
$ws.Range('E8').Value = '  -1.58%  '

# Row 9
$ws.Range('D9').Value = "'0.07857"
$ws.Range('E9').Value = '  -0.67%  '

# Row 10
$ws.Range('D10').Value = "'0.9566"
$ws.Range('E10').Value = '  -2.77%  '

# Row 11
$ws.Range('D11').Value = "'21.81"
$ws.Range('E11').Value = '  -1.48%  '

# Row 12
$ws.Range('D12').Value = "'1.885.92"
$ws.Range('E12').Value = '  -1.15%  '

# Row 13
$ws.Range('D13').Value = "'5.648"
$ws.Range('E13').Value = '  -3.39%  '

# Row 14
$ws.Range('D14').Value = "'6.877"
$ws.Range('E14').Value = '  -2.00%  '

# Row 15
$ws.Range('E15').Value = '  -0.65%  '

# Row 16
$ws.Range('D16').Value = "'1.001"
$ws.Range('E16').Value = '  -0.35%  '

# Row 17
$ws.Range('D17').Value = "'86.76"
$ws.Range('E17').Value = '  -0.99%  '

# Row 18
$ws.Range('D18').Value = "'0.000009905"
$ws.Range('E18').Value = '  -2.23%  '

# Row 19
$ws.Range('D19').Value = "'16.57"
$ws.Range('E19').Value = '  -2.75%  '

# Row 20
$ws.Range('E20').Value = '  -0.41%  '

# Row 21
$ws.Range('D21').Value = "'28.024.24"
$ws.Range('E21').Value = '  -2.19%  '

# Row 22
$ws.Range('D22').Value = "'5.295"
$ws.Range('E22').Value = '  -2.13%  '

# Row 23
$ws.Range('D23').Value = "'10.95"
$ws.Range('E23').Value = '  -3.05%  '

# Row 24
$ws.Range('D24').Value = "'2.087"
$ws.Range('E24').Value = '  -1.96%  '

# Row 25
$ws.Range('D25').Value = "'2.111.28"
$ws.Range('E25').Value = '  -0.50%  '

# Row 26
$ws.Range('D26').Value = "'153.80"
$ws.Range('E26').Value = '  +0.30%  '

# Row 27
$ws.Range('D27').Value = "'19.11"
$ws.Range('E27').Value = '  -1.79%  '

# Row 28
$ws.Range('D28').Value = "'5.716"
$ws.Range('E28').Value = '  -9.00%  '

# Row 29
$ws.Range('E29').Value = '  -3.01%  '

# Row 30
$ws.Range('D30').Value = "'116.95"
$ws.Range('E30').Value = '  -0.52%  '

# Row 31
$ws.Range('D31').Value = "'0.9345"
$ws.Range('E31').Value = '  -4.76%  '

# Row 32
$ws.Range('E32').Value = '  -2.17%  '

# Row 33
$ws.Range('D33').Value = "'5.285"
$ws.Range('E33').Value = '  -1.96%  '

# Row 34
$ws.Range('E34').Value = '  -2.66%  '

# Row 35
$ws.Range('D35').Value = "'3.319"
$ws.Range('E35').Value = '  -5.41%  '

# Row 36
$ws.Range('D36').Value = "'0.05854"
$ws.Range('E36').Value = '  -4.61%  '

# Row 37
$ws.Range('E37').Value = '  -2.66%  '

# Row 38
$ws.Range('D38').Value = "'1.144"
$ws.Range('E38').Value = '  -0.98%  '

# Row 39
$ws.Range('D39').Value = "'7.756"
$ws.Range('E39').Value = '  +1.89%  '

# Row 40
$ws.Range('D40').Value = "'0.5565"
$ws.Range('E40').Value = '  -2.71%  '

# Row 41
$ws.Range('D41').Value = "'9.865"
$ws.Range('E41').Value = '  -2.32%  '

# Row 42
$ws.Range('E42').Value = '  -1.68%  '

# Row 43
$ws.Range('E43').Value = '  -1.94%  '

# Row 44
$ws.Range('D44').Value = "'11.55"
$ws.Range('E44').Value = '  -2.81%  '

# Row 45
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = "'0.07020"
$ws.Range('E45').Value = '  -1.62%  '

# Row 46
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.5252"
$ws.Range('E46').Value = '  -2.74%  '

# Row 47
$ws.Range('D47').Value = "'2.136"
$ws.Range('E47').Value = '  -9.60%  '

# Row 48
$ws.Range('D48').Value = "'1.825"
$ws.Range('E48').Value = '  -4.41%  '

# Row 49
$ws.Range('D49').Value = "'112.56"
$ws.Range('E49').Value = '  -2.29%  '

# Row 50
$ws.Range('E50').Value = '  -0.21%  '

# Row 51
$ws.Range('D51').Value = "'2.319"
$ws.Range('E51').Value = '  +0.06%  '
